# Weekly update: a new price record was reported for the week, so a new
# row is inserted at row 16 (pushing the existing rows 16-70 down to 17-71)
# and populated with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("16:16").Insert()

$ws.Range("A16").Value = 11
$ws.Range("B16").Value = "Vega Monumental Concepción"
$ws.Range("C16").Value = "Bíobío"
$ws.Range("D16").Value = 44811
$ws.Range("E16").Value = 8
$ws.Range("F16").Value = 100112013
$ws.Range("G16").Value = "Alcachofa"
$ws.Range("H16").Value = "Madrigal"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 12000
$ws.Range("L16").Value = 13000
$ws.Range("M16").Value = 12500
$ws.Range("N16").Value = "$/caja 40 unidades"
$ws.Range("O16").Value = "Provincia de Limarí"
$ws.Range("P16").Value = 312
$ws.Range("Q16").Value = 40
$ws.Range("R16").Value = "Hortaliza"
